# Update Argent (silver/solar component) prices: append a new row (row 3)
# with date 2025-03-04 and the latest price to each of the 9 price sheets.
# Values are written as text (matching the existing inlineStr-style cells
# in row 2), so NumberFormat is forced to "@" before the assignment to stop
# Excel from auto-coercing these numeric/date-looking strings into real
# numbers or dates.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-04"

# Worksheet tab order matches the per-sheet price to append in column B.
$updates = @(
    @{ Sheet = 1; Price = "5.48" },    # N-Dense
    @{ Sheet = 2; Price = "5.89" },    # N-Type
    @{ Sheet = 3; Price = "1.19" },    # N-type Wafer
    @{ Sheet = 4; Price = "0.29" },    # Cell Topcon 183mm
    @{ Sheet = 5; Price = "0.1" },     # Module Topcon 183mm
    @{ Sheet = 6; Price = "5,179" },   # Silver Rear_side
    @{ Sheet = 7; Price = "7,753" },   # Silver Busbar front-side
    @{ Sheet = 8; Price = "7,803" },   # Silver finger front-side
    @{ Sheet = 9; Price = "7.2998" }   # USD_CNY
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    $dateCell = $ws.Cells.Item(3, 1)
    $priceCell = $ws.Cells.Item(3, 2)

    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $u.Price
}
